# update upload contract and add contractlist module testcases

$wb = $excel.ActiveWorkbook

# --- 1. Extend the "upload_contract" sheet with two more test-case rows ---
$upload = $wb.Worksheets.Item("upload_contract")

$upload.Range("A8").Value = "bss"

$upload.Range("A9").Value = "哈哈哈"
$upload.Range("B9").Value = "没有数据哦~"
$upload.Range("C9").Value = "搜索文档无数据"

# move the cursor to the last entered cell, matching the saved selection
$upload.Range("C9").Select()

# --- 2. Add the new "contract_list" module worksheet after "upload_contract" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$contractList = $wb.Worksheets.Add($null, $lastSheet)
$contractList.Name = "contract_list"

$contractList.Range("A1").Value = "reason"
$contractList.Range("B1").Value = "except_result"
$contractList.Range("C1").Value = "describe"

$contractList.Range("A2").Value = "拒绝契约，拒绝契约"
$contractList.Range("B2").Value = "操作成功"

$contractList.Range("B3").Value = "拒绝原因不能为空"
